$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# D1 header - "Post Treatment" (already exists as shared string index 2, same label as C1's sibling concept)
$ws.Range("D1").Value = "Post Treatment"

# Individual symptom values (rows 2-23) for column D ("Post Treatment")
$dValues = @{
    2  = 4
    3  = 4
    4  = 2
    5  = 5
    6  = 3
    7  = 4
    8  = 2
    9  = 4
    10 = 2
    11 = 2
    12 = 5
    13 = 3
    14 = 3
    15 = 3
    16 = 2
    17 = 2
    18 = 5
    19 = 6
    20 = 5
    21 = 3
    22 = 6
    23 = 2
}

foreach ($row in $dValues.Keys) {
    $ws.Cells.Item($row, 4).Value = $dValues[$row]
}

# Summary / cluster formulas for column D (rows 24-30), mirroring column C's formulas
$ws.Range("D24").Formula = "=SUM(D2:D23)"
$ws.Range("D25").Formula = "=SUM(D2 + D12 + D13)"
$ws.Range("D26").Formula = "=SUM(D19:D22)"
$ws.Range("D27").Formula = "=SUM(D14:D18)"
$ws.Range("D28").Formula = "=D23"
$ws.Range("D29").Formula = "=SUM(D3:D6)"
$ws.Range("D30").Formula = "=SUM(D7:D11)"

# Update the selection/active cell to match the recorded cursor position (row 21 selected)
$ws.Rows.Item(21).Select()
